$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Duplicate row 14 (the "OTHER" charge row) down into a new row 15,
# preserving all of row 14's cell formatting/styles.
$ws.Rows("14:14").Copy()
$ws.Rows("15:15").Insert()

# Turn (the now-duplicated) row 14 into the new "FURNITURES / DECORATION"
# charge entry for France.
$ws.Range("A14").Value = "FURNITURES / DECORATION"
$ws.Range("B14").Value = "FRANCE"
$ws.Range("C14").Value = "'/FRA"

# Row 15 keeps the original "OTHER" entry it inherited from the copy, but
# it never had the (empty, formatted-only) D column cell that row 14 had,
# so drop it.
$ws.Range("D15").Clear()
$ws.Rows("15:15").RowHeight = 15

$null = $ws.Range("A17").Select()
